$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.340.61'
$ws.Range('E2').Value = '  -4.27%  '

# Row 3
$ws.Range('D3').Value = '2.505.43'
$ws.Range('E3').Value = '  -5.35%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.85'
$ws.Range('E5').Value = '  -2.12%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.18'
$ws.Range('E6').Value = '  -4.82%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('E8').Value = '  -0.99%  '

# Row 9
$ws.Range('D9').Value = '2.502.61'
$ws.Range('E9').Value = '  -5.41%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.156'
$ws.Range('E10').Value = '  -9.31%  '

# Row 11
$ws.Range('E11').Value = '  -1.32%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.343'
$ws.Range('E12').Value = '  -3.68%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.84'
$ws.Range('E13').Value = '  -2.18%  '

# Row 14
$ws.Range('D14').Value = '2.963.33'
$ws.Range('E14').Value = '  -5.35%  '

# Row 15
$ws.Range('D15').Value = '69.281.67'
$ws.Range('E15').Value = '  -4.23%  '

# Row 16
$ws.Range('E16').Value = '  -7.07%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.70'
$ws.Range('E17').Value = '  -4.76%  '

# Row 18
$ws.Range('D18').Value = '2.499.87'
$ws.Range('E18').Value = '  -4.99%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  -6.97%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.19'
$ws.Range('E21').Value = '  -7.41%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.92'
$ws.Range('E22').Value = '  -5.62%  '

# Row 23
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.94'
$ws.Range('E23').Value = '  -6.59%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.19'
$ws.Range('E25').Value = '  -3.90%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.94'
$ws.Range('E26').Value = '  -7.25%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('E27').Value = '  -8.11%  '

# Row 28
$ws.Range('D28').Value = '2.632.12'
$ws.Range('E28').Value = '  -5.37%  '

# Row 29
$ws.Range('E29').Value = '  +0.10%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0893'
$ws.Range('E30').Value = '  -6.67%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.84'
$ws.Range('E31').Value = '  -2.40%  '

# Row 32
$ws.Range('E32').Value = '  -2.13%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '455.75'
$ws.Range('E33').Value = '  -8.47%  '

# Row 34
$ws.Range('E34').Value = '  -3.51%  '

# Row 35
$ws.Range('E35').Value = '  +0.00%  '

# Row 36
$ws.Range('E36').Value = '  +0.84%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '154.43'
$ws.Range('E37').Value = '  -5.51%  '

# Row 38
$ws.Range('E38').Value = '  +0.38%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.34'
$ws.Range('E39').Value = '  -4.90%  '

# Row 40
$ws.Range('E40').Value = '  +0.03%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.71'
$ws.Range('E41').Value = '  -4.06%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.315'
$ws.Range('E42').Value = '  -3.80%  '

# Row 43
$ws.Range('E43').Value = '  -8.86%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.16'
$ws.Range('E44').Value = '  -14.92%  '

# Row 45
$ws.Range('E45').Value = '  -2.60%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.28'
$ws.Range('E46').Value = '  -11.83%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.62'
$ws.Range('E47').Value = '  -6.68%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.525'
$ws.Range('E48').Value = '  -4.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.49'
$ws.Range('E49').Value = '  -4.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  -5.52%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0730'
$ws.Range('E51').Value = '  -2.33%  '
